$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-02-05 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-02-06 Thursday", 2) | Out-Null
$d.Content.Find.Execute("73÷3=24, 1", $true, $false, $false, $false, $false, $true, 1, $false, "95÷7=13, 4", 2) | Out-Null
$d.Content.Find.Execute("21÷7=3, 0", $true, $false, $false, $false, $false, $true, 1, $false, "93÷5=18, 3", 2) | Out-Null
$d.Content.Find.Execute("76÷2=38, 0", $true, $false, $false, $false, $false, $true, 1, $false, "48÷4=12, 0", 2) | Out-Null
$d.Content.Find.Execute("85÷4=21, 1", $true, $false, $false, $false, $false, $true, 1, $false, "42÷2=21, 0", 2) | Out-Null
$d.Content.Find.Execute("26÷5=5, 1", $true, $false, $false, $false, $false, $true, 1, $false, "99÷2=49, 1", 2) | Out-Null
$d.Content.Find.Execute("31÷3=10, 1", $true, $false, $false, $false, $false, $true, 1, $false, "70÷2=35, 0", 2) | Out-Null
$d.Content.Find.Execute("88÷5=17, 3", $true, $false, $false, $false, $false, $true, 1, $false, "69÷2=34, 1", 2) | Out-Null
$d.Content.Find.Execute("23÷4=5, 3", $true, $false, $false, $false, $false, $true, 1, $false, "93÷8=11, 5", 2) | Out-Null
$d.Content.Find.Execute("18÷7=2, 4", $true, $false, $false, $false, $false, $true, 1, $false, "69÷2=34, 1", 2) | Out-Null
$d.Content.Find.Execute("44÷6=7, 2", $true, $false, $false, $false, $false, $true, 1, $false, "72÷5=14, 2", 2) | Out-Null
$d.Content.Find.Execute("25÷6=4, 1", $true, $false, $false, $false, $false, $true, 1, $false, "68÷4=17, 0", 2) | Out-Null
$d.Content.Find.Execute("46÷8=5, 6", $true, $false, $false, $false, $false, $true, 1, $false, "39÷7=5, 4", 2) | Out-Null
$d.Content.Find.Execute("38÷9=4, 2", $true, $false, $false, $false, $false, $true, 1, $false, "83÷4=20, 3", 2) | Out-Null
$d.Content.Find.Execute("79÷3=26, 1", $true, $false, $false, $false, $false, $true, 1, $false, "35÷8=4, 3", 2) | Out-Null
$d.Content.Find.Execute("40÷2=20, 0", $true, $false, $false, $false, $false, $true, 1, $false, "98÷3=32, 2", 2) | Out-Null
$d.Content.Find.Execute("64÷8=8, 0", $true, $false, $false, $false, $false, $true, 1, $false, "88÷9=9, 7", 2) | Out-Null
$d.Content.Find.Execute("18÷4=4, 2", $true, $false, $false, $false, $false, $true, 1, $false, "23÷6=3, 5", 2) | Out-Null
$d.Content.Find.Execute("97÷8=12, 1", $true, $false, $false, $false, $false, $true, 1, $false, "97÷9=10, 7", 2) | Out-Null
$d.Content.Find.Execute("75÷7=10, 5", $true, $false, $false, $false, $false, $true, 1, $false, "13÷6=2, 1", 2) | Out-Null
$d.Content.Find.Execute("88÷7=12, 4", $true, $false, $false, $false, $false, $true, 1, $false, "93÷8=11, 5", 2) | Out-Null
$d.Content.Find.Execute("59÷5=11, 4", $true, $false, $false, $false, $false, $true, 1, $false, "98÷3=32, 2", 2) | Out-Null
$d.Content.Find.Execute("13÷2=6, 1", $true, $false, $false, $false, $false, $true, 1, $false, "34÷9=3, 7", 2) | Out-Null
$d.Content.Find.Execute("98÷5=19, 3", $true, $false, $false, $false, $false, $true, 1, $false, "79÷7=11, 2", 2) | Out-Null
$d.Content.Find.Execute("43÷7=6, 1", $true, $false, $false, $false, $false, $true, 1, $false, "68÷9=7, 5", 2) | Out-Null
$d.Content.Find.Execute("87÷7=12, 3", $true, $false, $false, $false, $false, $true, 1, $false, "42÷4=10, 2", 2) | Out-Null
